$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts the old row 8 and everything below it down by one).
# This new row 8 becomes the "DK/NA" entry belonging to the "Interested in News" category.
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = ""
$ws.Range("B8").Value = "DK/NA"
$ws.Range("C8").Value = "(0.0%) 0"
$ws.Range("D8").Value = "(0.0%) 0"

# Append a new "DK/NA" row after the last existing row (old row 13 "Woman", now row 14)
# for the "Gender" category, i.e. the new row 15.
$ws.Range("A15").Value = ""
$ws.Range("B15").Value = "DK/NA"
$ws.Range("C15").Value = "(0.0%) 0"
$ws.Range("D15").Value = "(0.0%) 0"
